$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 236, shifting the
# existing rows 236-266 down to 238-268 (weekly update adds two new
# observations for this market/category combination).
$ws.Range("A236:A237").EntireRow.Insert()

# New row 236: Primera quality reading for the new date.
$ws.Cells.Item(236, 1).Value = 9
$ws.Cells.Item(236, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(236, 3).Value = "Metropolitana"
$ws.Cells.Item(236, 4).Value = 44491
$ws.Cells.Item(236, 5).Value = 13
$ws.Cells.Item(236, 6).Value = 100112012
$ws.Cells.Item(236, 7).Value = "Espinaca"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 250
$ws.Cells.Item(236, 11).Value = 5000
$ws.Cells.Item(236, 12).Value = 6000
$ws.Cells.Item(236, 13).Value = 5500
$ws.Cells.Item(236, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(236, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(236, 16).Value = 550
$ws.Cells.Item(236, 17).Value = 10
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# New row 237: Segunda quality reading for the same new date.
$ws.Cells.Item(237, 1).Value = 9
$ws.Cells.Item(237, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(237, 3).Value = "Metropolitana"
$ws.Cells.Item(237, 4).Value = 44491
$ws.Cells.Item(237, 5).Value = 13
$ws.Cells.Item(237, 6).Value = 100112012
$ws.Cells.Item(237, 7).Value = "Espinaca"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Segunda"
$ws.Cells.Item(237, 10).Value = 97
$ws.Cells.Item(237, 11).Value = 4000
$ws.Cells.Item(237, 12).Value = 4500
$ws.Cells.Item(237, 13).Value = 4253
$ws.Cells.Item(237, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(237, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(237, 16).Value = 425
$ws.Cells.Item(237, 17).Value = 10
$ws.Cells.Item(237, 18).Value = "Hortaliza"
